# Demographic.xlsx edit:
#  - B2 header text corrected from "group" to "Group"
#  - Row 5 (A5, B5) values cleared, keeping the existing cell formatting
#  - Active selection moved to I10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix capitalization of the "group" label in B2 -> "Group"
$ws.Range("B2").Value = "Group"

# Clear the now-unused row 5 data (was A5=4, B5="Group")
$ws.Range("A5:B5").ClearContents()

# Move/restore the active cell selection to I10
$ws.Range("I10").Select()
